$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.306.86'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.481.57'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.20%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.480.95'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.22%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.143'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.56'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.429'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('E13').Value = '  -2.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '31.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.069.96'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.482.93'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.317.61'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('E19').Value = '  +1.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '445.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.625'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.86%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.624.49'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.05%  '
$ws.Range('E27').Value = '  -4.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.94'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.50'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('E32').Value = '  +4.19%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.47'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.476.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  +6.58%  '
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '177.24'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.890'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '30.10'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.99%  '
$ws.Range('E49').Value = '  -4.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('E51').Value = '  -1.36%  '
